# Semana 50 de 2025: add week-50 column (BA) to the weekly IRA report,
# and correct a handful of previously-entered weekly counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header label for week 50 (text, like the other week-number headers) ---
$ws.Cells.Item(1, 53).Value = "'50"

# --- Per-UPGD weekly counts: new week-50 values (column BA) plus a few corrections ---
$ws.Cells.Item(2, 53).Value = 20  # BA2
$ws.Cells.Item(3, 53).Value = 58  # BA3
$ws.Cells.Item(5, 52).Value = 1  # AZ5
$ws.Cells.Item(5, 53).Value = 3  # BA5
$ws.Cells.Item(6, 53).Value = 77  # BA6
$ws.Cells.Item(7, 53).Value = 27  # BA7
$ws.Cells.Item(8, 53).Value = 23  # BA8
$ws.Cells.Item(12, 53).Value = 4  # BA12
$ws.Cells.Item(13, 53).Value = 1  # BA13
$ws.Cells.Item(14, 53).Value = 1  # BA14
$ws.Cells.Item(15, 53).Value = 1  # BA15
$ws.Cells.Item(16, 53).Value = 2  # BA16
$ws.Cells.Item(17, 53).Value = 1  # BA17
$ws.Cells.Item(19, 53).Value = 1  # BA19
$ws.Cells.Item(23, 53).Value = 5  # BA23
$ws.Cells.Item(25, 53).Value = 44  # BA25
$ws.Cells.Item(28, 53).Value = 208  # BA28
$ws.Cells.Item(29, 53).Value = 0  # BA29
$ws.Cells.Item(30, 52).Value = 33  # AZ30
$ws.Cells.Item(30, 53).Value = 51  # BA30
$ws.Cells.Item(31, 53).Value = 2  # BA31
$ws.Cells.Item(35, 53).Value = 53  # BA35
$ws.Cells.Item(36, 53).Value = 3  # BA36
$ws.Cells.Item(38, 53).Value = 66  # BA38
$ws.Cells.Item(41, 53).Value = 5  # BA41
$ws.Cells.Item(42, 53).Value = 79  # BA42
$ws.Cells.Item(43, 29).Value = 54  # AC43
$ws.Cells.Item(43, 35).Value = 55  # AI43
$ws.Cells.Item(43, 36).Value = 54  # AJ43
$ws.Cells.Item(43, 37).Value = 41  # AK43
$ws.Cells.Item(43, 38).Value = 55  # AL43
$ws.Cells.Item(43, 39).Value = 53  # AM43
$ws.Cells.Item(43, 40).Value = 74  # AN43
$ws.Cells.Item(43, 41).Value = 55  # AO43
$ws.Cells.Item(43, 42).Value = 50  # AP43
$ws.Cells.Item(43, 43).Value = 46  # AQ43
$ws.Cells.Item(43, 44).Value = 42  # AR43
$ws.Cells.Item(43, 45).Value = 45  # AS43
$ws.Cells.Item(43, 46).Value = 50  # AT43
$ws.Cells.Item(43, 47).Value = 42  # AU43
$ws.Cells.Item(43, 48).Value = 35  # AV43
$ws.Cells.Item(43, 53).Value = 45  # BA43
$ws.Cells.Item(46, 53).Value = 86  # BA46
$ws.Cells.Item(47, 53).Value = 113  # BA47
$ws.Cells.Item(48, 53).Value = 3  # BA48
$ws.Cells.Item(49, 53).Value = 100  # BA49
$ws.Cells.Item(50, 53).Value = 3  # BA50
$ws.Cells.Item(51, 53).Value = 0  # BA51
$ws.Cells.Item(53, 52).Value = 3  # AZ53
$ws.Cells.Item(53, 53).Value = 4  # BA53
$ws.Cells.Item(54, 53).Value = 10  # BA54
$ws.Cells.Item(55, 53).Value = 1  # BA55
$ws.Cells.Item(56, 53).Value = 0  # BA56
$ws.Cells.Item(57, 13).Value = 9  # M57
$ws.Cells.Item(57, 14).Value = 5  # N57
$ws.Cells.Item(57, 32).Value = 8  # AF57
$ws.Cells.Item(57, 33).Value = 16  # AG57
$ws.Cells.Item(57, 34).Value = 5  # AH57
$ws.Cells.Item(57, 53).Value = 11  # BA57
$ws.Cells.Item(58, 53).Value = 35  # BA58
$ws.Cells.Item(59, 53).Value = 2  # BA59
